$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# evidence_extraction_quality: 3 -> 4
$ws.Range("B2").Value = 4

# w_evidence_extraction_quality: 0.75 -> 1
$ws.Range("H2").Value = 1

# weighted_final_score: 3.3 -> 3.55
$ws.Range("N2").Value = 3.55

# justification text update
$ws.Range("O2").Value = "The report demonstrates strong evidence extraction quality with accurate citations and full sentence quotations, though it lacks some depth in coverage of representativeness dimensions, particularly in demographic and geographic specifics. The structure is clear and well-organized, aiding readability. The relevance and faithfulness of the evidence are maintained, with no unsupported assumptions. However, the identification of missing disclosures could be more detailed, particularly regarding specific demographic and geographic targets. The audit usefulness is reasonable, but the report could benefit from more explicit traceability and verifiability of claims. Overall, the report is strong but has clear areas for improvement, especially in detailing missing disclosures and enhancing audit traceability."
